$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.399.01'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.847.09'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6311'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07554'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2953'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('E10').Value = '  -0.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07691'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.860.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.987'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6839'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.05'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.108.74'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.124'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.431.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '228.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9997'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.537'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.69%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1392'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.69%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.371'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.28%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.67'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.470'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05731'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.265'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.122'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.36%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.022'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.844'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.08%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.155'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.28%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7162'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.78%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.588'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.250.71'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01807'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.777'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9092'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.175'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.74%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.000'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.60'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.08'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.69%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.099'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.41%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000118'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4015'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.090'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.678'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.13%  '
$ws.Range('E51').Value = '  +0.00%  '
